$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G3").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G4").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G5").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G6").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G8").Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G10").Value = "Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G13").Value = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, maimustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G14").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G15").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G16").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G17").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G21").Value = "Alshimaa_khaled@med.asu.edu.eg, System"
$ws.Range("G23").Value = "menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G24").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G25").Value = "youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G27").Value = "nourhanmohamed@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G29").Value = "nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G33").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G34").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G35").Value = "gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G36").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G37").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G40").Value = "Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G41").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G43").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G44").Value = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, maimustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G45").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G46").Value = "Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G48").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G52").Value = "afnan.fares@med.asu.edu.eg, System"
$ws.Range("G54").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G55").Value = "menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G56").Value = "youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G58").Value = "nourhanmohamed@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G60").Value = "nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G64").Value = "mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G65").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G66").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G67").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G68").Value = "rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G70").Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G71").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G72").Value = "Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G73").Value = "Safa.hany@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg"
$ws.Range("G74").Value = "Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg"
$ws.Range("G75").Value = "Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G76").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G77").Value = "Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G79").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G83").Value = "Alshimaa_khaled@med.asu.edu.eg, System"
$ws.Range("G85").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G86").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G89").Value = "nourhan.mostafa@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G91").Value = "nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G95").Value = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G96").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G97").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G98").Value = "eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G99").Value = "nesmadrahim@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G101").Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G102").Value = "Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G103").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G105").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G106").Value = "nadia.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G107").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G108").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G109").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G114").Value = "afnan.fares@med.asu.edu.eg, System"
$ws.Range("G116").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg"
$ws.Range("G117").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg"
$ws.Range("G120").Value = "nourhanmohamed@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G122").Value = "neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G126").Value = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G127").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G128").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G129").Value = "eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G130").Value = "nesmadrahim@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
$ws.Range("G132").Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G133").Value = "Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G136").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G137").Value = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, maimustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G138").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G139").Value = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G140").Value = "Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G141").Value = "Alshimaa_khaled@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Range("G145").Value = "Alshimaa_khaled@med.asu.edu.eg, System"
$ws.Range("G147").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg"
$ws.Range("G148").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg"
$ws.Range("G153").Value = "neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G154").Value = "eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G157").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G158").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G159").Value = "gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G160").Value = "servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G161").Value = "rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G164").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G166").Value = "Safa.hany@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg"
$ws.Range("G167").Value = "Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, mernamahrous@med.asu.edu.eg"
$ws.Range("G168").Value = "Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G169").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G170").Value = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G172").Value = "Alshimaa_khaled@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Range("G176").Value = "Alshimaa_khaled@med.asu.edu.eg, System"
$ws.Range("G178").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G179").Value = "menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G182").Value = "nourhanmohamed@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G184").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G185").Value = "eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
